$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 7.861494
$ws.Range("H2").Value = 23.584482
$ws.Range("I2").Value = 0.02087137335661869
$ws.Range("J2").Value = 0.02087137335661869
$ws.Range("M2").Value = 8.495336333333334
$ws.Range("N2").Value = 25.486009
$ws.Range("O2").Value = 0.9246706102479199
$ws.Range("P2").Value = 0.9246706102479199
$ws.Range("Q2").Value = 66.78603561248201
$ws.Range("R2").Value = 601.074320512338
$ws.Range("S2").Value = 0.01929914553837678
$ws.Range("T2").Value = 0.01929914553837678
$ws.Range("G3").Value = 7.861494
$ws.Range("H3").Value = 23.584482
$ws.Range("I3").Value = 0.02087137335661869
$ws.Range("J3").Value = 0.02087137335661869
$ws.Range("O3").Value = 0.06049852883963747
$ws.Range("P3").Value = 0.06049852883963747
$ws.Range("Q3").Value = 4.369617523047999
$ws.Range("R3").Value = 39.32655770743199
$ws.Range("S3").Value = 0.001262687382938237
$ws.Range("T3").Value = 0.001262687382938237
$ws.Range("G4").Value = 7.861494
$ws.Range("H4").Value = 23.584482
$ws.Range("I4").Value = 0.02087137335661869
$ws.Range("J4").Value = 0.02087137335661869
$ws.Range("M4").Value = 0.1362573333333333
$ws.Range("N4").Value = 0.408772
$ws.Range("O4").Value = 0.01483086091244269
$ws.Range("P4").Value = 0.01483086091244269
$ws.Range("Q4").Value = 1.071186208456
$ws.Range("R4").Value = 9.640675876103998
$ws.Range("S4").Value = 0.0003095404353036739
$ws.Range("T4").Value = 0.0003095404353036739
$ws.Range("I5").Value = 0.9496861641109521
$ws.Range("J5").Value = 0.9496861641109521
$ws.Range("M5").Value = 8.495336333333334
$ws.Range("N5").Value = 25.486009
$ws.Range("O5").Value = 0.9246706102479199
$ws.Range("P5").Value = 0.9246706102479199
$ws.Range("Q5").Value = 3038.888380427636
$ws.Range("R5").Value = 27349.99542384872
$ws.Range("S5").Value = 0.8781468849124803
$ws.Range("T5").Value = 0.8781468849124803
$ws.Range("I6").Value = 0.9496861641109521
$ws.Range("J6").Value = 0.9496861641109521
$ws.Range("O6").Value = 0.06049852883963747
$ws.Range("P6").Value = 0.06049852883963747
$ws.Range("S6").Value = 0.05745461578807112
$ws.Range("T6").Value = 0.05745461578807112
$ws.Range("I7").Value = 0.9496861641109521
$ws.Range("J7").Value = 0.9496861641109521
$ws.Range("M7").Value = 0.1362573333333333
$ws.Range("N7").Value = 0.408772
$ws.Range("O7").Value = 0.01483086091244269
$ws.Range("P7").Value = 0.01483086091244269
$ws.Range("Q7").Value = 48.74095748158
$ws.Range("R7").Value = 438.66861733422
$ws.Range("S7").Value = 0.01408466341040075
$ws.Range("T7").Value = 0.01408466341040075
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 7.100387666666667
$ws.Range("H8").Value = 21.301163
$ws.Range("I8").Value = 0.01885072251759407
$ws.Range("J8").Value = 0.01885072251759407
$ws.Range("M8").Value = 8.495336333333334
$ws.Range("N8").Value = 25.486009
$ws.Range("O8").Value = 0.9246706102479199
$ws.Range("P8").Value = 0.9246706102479199
$ws.Range("Q8").Value = 60.32018132538524
$ws.Range("R8").Value = 542.8816319284671
$ws.Range("S8").Value = 0.01743070909395791
$ws.Range("T8").Value = 0.01743070909395791
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 7.100387666666667
$ws.Range("H9").Value = 21.301163
$ws.Range("I9").Value = 0.01885072251759407
$ws.Range("J9").Value = 0.01885072251759407
$ws.Range("O9").Value = 0.06049852883963747
$ws.Range("P9").Value = 0.06049852883963747
$ws.Range("Q9").Value = 3.946575341620889
$ws.Range("R9").Value = 35.51917807458801
$ws.Range("S9").Value = 0.001140440979878668
$ws.Range("T9").Value = 0.001140440979878668
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 7.100387666666667
$ws.Range("H10").Value = 21.301163
$ws.Range("I10").Value = 0.01885072251759407
$ws.Range("J10").Value = 0.01885072251759407
$ws.Range("M10").Value = 0.1362573333333333
$ws.Range("N10").Value = 0.408772
$ws.Range("O10").Value = 0.01483086091244269
$ws.Range("P10").Value = 0.01483086091244269
$ws.Range("Q10").Value = 0.9674798890928888
$ws.Range("R10").Value = 8.707319001836
$ws.Range("S10").Value = 0.0002795724437574891
$ws.Range("T10").Value = 0.0002795724437574891
$ws.Range("G11").Value = 2.906965
$ws.Range("H11").Value = 8.720895000000001
$ws.Range("I11").Value = 0.007717661789174303
$ws.Range("J11").Value = 0.007717661789174304
$ws.Range("M11").Value = 8.495336333333334
$ws.Range("N11").Value = 25.486009
$ws.Range("O11").Value = 0.9246706102479199
$ws.Range("P11").Value = 0.9246706102479199
$ws.Range("Q11").Value = 24.69564538422834
$ws.Range("R11").Value = 222.260808458055
$ws.Range("S11").Value = 0.007136295036282856
$ws.Range("T11").Value = 0.007136295036282857
$ws.Range("G12").Value = 2.906965
$ws.Range("H12").Value = 8.720895000000001
$ws.Range("I12").Value = 0.007717661789174303
$ws.Range("J12").Value = 0.007717661789174304
$ws.Range("O12").Value = 0.06049852883963747
$ws.Range("P12").Value = 0.06049852883963747
$ws.Range("Q12").Value = 1.615764790113333
$ws.Range("R12").Value = 14.54188311102
$ws.Range("S12").Value = 0.0004669071843269297
$ws.Range("T12").Value = 0.0004669071843269298
$ws.Range("G13").Value = 2.906965
$ws.Range("H13").Value = 8.720895000000001
$ws.Range("I13").Value = 0.007717661789174303
$ws.Range("J13").Value = 0.007717661789174304
$ws.Range("M13").Value = 0.1362573333333333
$ws.Range("N13").Value = 0.408772
$ws.Range("O13").Value = 0.01483086091244269
$ws.Range("P13").Value = 0.01483086091244269
$ws.Range("Q13").Value = 0.3960952989933333
$ws.Range("R13").Value = 3.56485769094
$ws.Range("S13").Value = 0.0001144595685645177
$ws.Range("T13").Value = 0.0001144595685645177
$ws.Range("G14").Value = 1.082561666666667
$ws.Range("H14").Value = 3.247685
$ws.Range("I14").Value = 0.002874078225660846
$ws.Range("J14").Value = 0.002874078225660847
$ws.Range("M14").Value = 8.495336333333334
$ws.Range("N14").Value = 25.486009
$ws.Range("O14").Value = 0.9246706102479199
$ws.Range("P14").Value = 0.9246706102479199
$ws.Range("Q14").Value = 9.196725459907222
$ws.Range("R14").Value = 82.77052913916501
$ws.Range("S14").Value = 0.002657575666822073
$ws.Range("T14").Value = 0.002657575666822074
$ws.Range("G15").Value = 1.082561666666667
$ws.Range("H15").Value = 3.247685
$ws.Range("I15").Value = 0.002874078225660846
$ws.Range("J15").Value = 0.002874078225660847
$ws.Range("O15").Value = 0.06049852883963747
$ws.Range("P15").Value = 0.06049852883963747
$ws.Range("Q15").Value = 0.6017151992288888
$ws.Range("R15").Value = 5.41543679306
$ws.Range("S15").Value = 0.0001738775044225168
$ws.Range("T15").Value = 0.0001738775044225168
$ws.Range("G16").Value = 1.082561666666667
$ws.Range("H16").Value = 3.247685
$ws.Range("I16").Value = 0.002874078225660846
$ws.Range("J16").Value = 0.002874078225660847
$ws.Range("M16").Value = 0.1362573333333333
$ws.Range("N16").Value = 0.408772
$ws.Range("O16").Value = 0.01483086091244269
$ws.Range("P16").Value = 0.01483086091244269
$ws.Range("Q16").Value = 0.1475069658688889
$ws.Range("R16").Value = 1.32756269282
$ws.Range("S16").Value = 0.00004262505441625609
$ws.Range("T16").Value = 0.00004262505441625609
